{"js": "// Bump the document revision number (7 -> 8) in the title-block paragraph,\n// split the \"January 2, 2020\" date into \"January 20, 2020\", and bump the\n// matching revision number (7 -> 8) in the primary footer.\n\n// --- 1 & 2: body paragraph \"Document version: 2020.1-EN-rev7, / based on\n// Curriculum - Version V5.1-EN; January 2, 2020\" ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet versionParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Document version:\") !== -1) {\n    versionParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (versionParagraph) {\n  // \"...rev7,\" -> \"...rev8,\"  (only the lone \"7\" run inside this paragraph)\n  const revResults = versionParagraph.search(\"7\", { matchCase: true });\n  await context.sync();\n  revResults.load(\"text\");\n  await context.sync();\n  if (revResults.items.length > 0) {\n    revResults.items[0].insertText(\"8\", \"Replace\");\n    await context.sync();\n  }\n\n  // \"January 2, 2020\" -> \"January 20, 2020\"\n  const dateResults = versionParagraph.search(\"January 2\", { matchCase: true });\n  await context.sync();\n  dateResults.load(\"text\");\n  await context.sync();\n  if (dateResults.items.length > 0) {\n    dateResults.items[0].insertText(\"0\", \"End\");\n    await context.sync();\n  }\n}\n\n// --- 3: primary footer \"Version 2020.1-EN-rev7\" -> \"...rev8\" ---\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nif (sections.items.length > 0) {\n  const primaryFooter = sections.items[0].getFooter(\"Primary\");\n  const footerRevResults = primaryFooter.search(\"7\", { matchCase: true });\n  await context.sync();\n  footerRevResults.load(\"text\");\n  await context.sync();\n  if (footerRevResults.items.length > 0) {\n    footerRevResults.items[0].insertText(\"8\", \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Bump the document revision number (7 -> 8) in the title-block paragraph,\n# split the \"January 2, 2020\" date into \"January 20, 2020\", and bump the\n# matching revision number (7 -> 8) in the primary footer.\n\n$d = $word.ActiveDocument\n\n# --- locate the title-block paragraph: \"Document version: 2020.1-EN-rev7,\n# based on Curriculum - Version V5.1-EN; January 2, 2020\" ---\n$versionParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($candidate.Range.Text -like \"*Document version:*\") {\n        $versionParagraph = $candidate\n        break\n    }\n}\n\nif ($versionParagraph -ne $null) {\n    # \"...rev7,\" -> \"...rev8,\"\n    $find1 = $versionParagraph.Range.Find\n    $find1.Execute(\"7\", $true, $false, $false, $false, $false, $true, 0, $false, \"8\", 1) | Out-Null\n\n    # \"January 2, 2020\" -> \"January 20, 2020\"\n    $find2 = $versionParagraph.Range.Find\n    $find2.Execute(\"January 2, 2020\", $true, $false, $false, $false, $false, $true, 0, $false, \"January 20, 2020\", 1) | Out-Null\n}\n\n# --- primary footer: \"Version 2020.1-EN-rev7\" -> \"...rev8\" ---\n$section = $d.Sections.Item(1)\n$footer = $section.Footers.Item(1)\n$find3 = $footer.Range.Find\n$find3.Execute(\"7\", $true, $false, $false, $false, $false, $true, 0, $false, \"8\", 1) | Out-Null\n"}
